# Updates cryptos list prices/volumes (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.754.49"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "2.908.99"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.23"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.62"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "2.906.62"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("E10").Value = "  -3.85%  "
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.91"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").Value = "3.391.00"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "61.833.44"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").Value = "2.908.18"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "435.51"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.657"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.94"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.83"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.14"
$ws.Range("E26").Value = "  -7.01%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000107"
$ws.Range("E29").Value = "  +21.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.976"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.04"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.51"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.07"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.35"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.271"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.57"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").Value = "2.697.78"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.29"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0336"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "343.82"
$ws.Range("E48").Value = "  -7.26%  "
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.24"
$ws.Range("E51").Value = "  -4.46%  "
